$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header date value for column BB, row 1 (matches style/format of existing BA1 date cell)
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("BB1").Value2 = 45986

# Rows 3-18: copy the same value as column BA (unchanged forecast continuation)
$sameRows = 3..18
foreach ($r in $sameRows) {
    $baVal = $ws.Cells.Item($r, 53).Value2
    $ws.Cells.Item($r, 54).Value2 = $baVal
}

# Rows 19-21: new distinct forecast values
$ws.Cells.Item(19, 54).Value2 = 2.46481303148316
$ws.Cells.Item(20, 54).Value2 = 3.633434696013671
$ws.Cells.Item(21, 54).Value2 = 3.559889218615653
